$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 647.5
$ws.Range("J18").Value = 995
$ws.Range("L18").Value = 995
$ws.Range("N18").Value = -1563
$ws.Range("H20").Value = 5867.25
$ws.Range("I20").Value = 5867.25
$ws.Range("K20").Value = 5867.25
$ws.Range("M20").Value = -5637.25
$ws.Range("H35").Value = 5867.25
$ws.Range("I35").Value = 5867.25
$ws.Range("K35").Value = 5867.25
$ws.Range("M35").Value = -5488.25
$ws.Range("H48").Value = 5000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 5000
$ws.Range("K48").Value = 0
$ws.Range("L48").ClearContents()
$ws.Range("M48").Value = 15000
$ws.Range("N48").Value = -15584
$ws.Range("H56").Value = 5000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 5000
$ws.Range("K56").Value = 0
$ws.Range("L56").ClearContents()
$ws.Range("M56").Value = 15000
$ws.Range("N56").Value = -16068
$ws.Range("H98").Value = 7910.778
$ws.Range("I98").Value = 2099.3333
$ws.Range("K98").Value = 2099.3333
$ws.Range("M98").Value = -601.3332999999998
$ws.Range("H122").Value = 7910.778
$ws.Range("I122").Value = 2099.3333
$ws.Range("K122").Value = 6297.999899999999
$ws.Range("M122").Value = -3847.999899999999
$ws.Range("H138").Value = 5226.8184
$ws.Range("J138").Value = 4199
$ws.Range("L138").Value = 12597
$ws.Range("N138").Value = -22877
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 38.857143
$ws.Range("I5").Value = 47.5
$ws.Range("K5").Value = 47.5
$ws.Range("M5").Value = 64.5
$ws.Range("H61").Value = 7457.1816
$ws.Range("J61").Value = 8447.5
$ws.Range("L61").Value = 8447.5
$ws.Range("N61").Value = -8871.5
$ws.Range("H74").Value = 4787.778
$ws.Range("J74").Value = 6583.3335
$ws.Range("L74").Value = 6583.3335
$ws.Range("N74").Value = -8331.333500000001
$ws.Range("H77").Value = 4787.778
$ws.Range("J77").Value = 6583.3335
$ws.Range("L77").Value = 32916.6675
$ws.Range("N77").Value = -41652.6675
$ws.Range("H97").Value = 2019.6666
$ws.Range("I97").Value = 2116.4167
$ws.Range("J97").Value = 1632.6666
$ws.Range("K97").Value = 2116.4167
$ws.Range("L97").Value = 1632.6666
$ws.Range("M97").Value = -1620.4167
$ws.Range("N97").Value = -2624.6666
$ws.Range("H102").Value = 1600.5
$ws.Range("I102").Value = 1600.5
$ws.Range("K102").Value = 1600.5
$ws.Range("M102").Value = 21.5
$ws.Range("H132").Value = 17666.666
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 17666.666
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").Value = 52999.99800000001
$ws.Range("N132").Value = -58059.99800000001
$ws.Range("H136").Value = 7457.1816
$ws.Range("J136").Value = 8447.5
$ws.Range("L136").Value = 25342.5
$ws.Range("N136").Value = -30442.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 38.857143
$ws.Range("I4").Value = 47.5
$ws.Range("K4").Value = 47.5
$ws.Range("M4").Value = 67.5
$ws.Range("H99").Value = 1074.875
$ws.Range("I99").Value = 1071.2858
$ws.Range("K99").Value = 1071.2858
$ws.Range("M99").Value = 426.7141999999999
$ws.Range("H105").Value = 5000
$ws.Range("I105").Value = 5000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5000
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -3253
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 159.05
$ws.Range("J7").Value = 205.44827
$ws.Range("L7").Value = 205.44827
$ws.Range("N7").Value = -431.44827
$ws.Range("H22").Value = 408.375
$ws.Range("I22").Value = 211.33333
$ws.Range("K22").Value = 211.33333
$ws.Range("M22").Value = 138.66667
$ws.Range("H31").Value = 5043.3335
$ws.Range("J31").Value = 6248.375
$ws.Range("L31").Value = 6248.375
$ws.Range("N31").Value = -6838.375
$ws.Range("H34").Value = 5043.3335
$ws.Range("J34").Value = 6248.375
$ws.Range("L34").Value = 6248.375
$ws.Range("N34").Value = -6652.375
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H44").Value = 14000
$ws.Range("J44").Value = 14000
$ws.Range("L44").Value = 14000
$ws.Range("N44").Value = -14884
$ws.Range("H122").Value = 1400
$ws.Range("I122").Value = 1400
$ws.Range("K122").Value = 4200
$ws.Range("M122").Value = -1750
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 166668830
$ws.Range("I4").Value = 2599.2
$ws.Range("K4").Value = 7797.599999999999
$ws.Range("M4").Value = -7685.599999999999
$ws.Range("H13").Value = 250
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 250
$ws.Range("K13").Value = 0
$ws.Range("L13").ClearContents()
$ws.Range("M13").Value = 750
$ws.Range("N13").Value = -1086
$ws.Range("H34").Value = 1296.4166
$ws.Range("J34").Value = 1933.375
$ws.Range("L34").Value = 5800.125
$ws.Range("N34").Value = -5968.125
$ws.Range("H50").Value = 265
$ws.Range("I50").Value = 265
$ws.Range("K50").Value = 795
$ws.Range("M50").Value = -314
$ws.Range("H53").Value = 265
$ws.Range("I53").Value = 265
$ws.Range("K53").Value = 795
$ws.Range("M53").Value = -314
$ws.Range("H131").Value = 11372.923
$ws.Range("I131").Value = 19299.428
$ws.Range("K131").Value = 57898.284
$ws.Range("M131").Value = -52858.284
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40518
$ws.Range("H102").Value = 3950
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 3950
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").Value = 3950
$ws.Range("N102").Value = -7194
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5750
$ws.Range("I7").Value = 5750
$ws.Range("K7").Value = 5750
$ws.Range("M7").Value = -5638
$ws.Range("H22").Value = 901
$ws.Range("I22").Value = 901
$ws.Range("K22").Value = 901
$ws.Range("M22").Value = -606
$ws.Range("H27").Value = 901
$ws.Range("I27").Value = 901
$ws.Range("K27").Value = 901
$ws.Range("M27").Value = -794
$ws.Range("H126").Value = 5750
$ws.Range("I126").Value = 5750
$ws.Range("K126").Value = 17250
$ws.Range("M126").Value = -14780
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 52975
$ws.Range("J62").Value = 52975
$ws.Range("L62").Value = 52975
$ws.Range("N62").Value = -54223
$ws.Range("H65").Value = 52975
$ws.Range("J65").Value = 52975
$ws.Range("L65").Value = 264875
$ws.Range("N65").Value = -271115
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 5785.278
$ws.Range("I132").Value = 3495.1667
$ws.Range("J132").Value = 10365.5
$ws.Range("K132").Value = 10485.5001
$ws.Range("L132").Value = 31096.5
$ws.Range("M132").Value = -7955.500100000001
$ws.Range("N132").Value = -36156.5
